# Add 2022-Q1 data:
#  - The existing "总计" sheet becomes the new "2022-Q1" detail sheet (21 funds).
#  - A brand new "总计" sheet (cloned from the original, so it keeps the same
#    sheetPr/format/structure) is appended with the updated totals table
#    (now including the 2022-Q1 row).

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")

# Clone the "总计" sheet *before* touching it, placing the clone immediately
# after it. The clone will become the refreshed "总计" sheet, the original
# will become "2022-Q1".
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item("总计 (2)")

$newTotal.Name = "总计_tmp"
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

$q1 = $oldTotal
$total = $newTotal

# ---------------------------------------------------------------------------
# Step 1: populate "2022-Q1" with the fund holdings table.
# ---------------------------------------------------------------------------

# Extend the header-row formatting (bold + border, same style as B1) across
# the new columns E1:H1, and extend the index-column formatting (same style
# as A2) down through A7:A22.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$q1.Range("A7:A22").PasteSpecial(-4122)

# Headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund holding rows (index, code, name, size, stock-position, position-pct, mv, rank)
$q1Rows = @(
    @(0, "006102", "浙商丰利增强债券",           "48.75", "47.92", "2.04", "0.9945", 8),
    @(1, "005535", "泰信竞争优选灵活配置混合",     "10.07", "90.21", "4.34", "0.4370", 4),
    @(2, "290006", "泰信蓝筹精选混合",             "7.62",  "88.52", "4.30", "0.3277", 4),
    @(3, "630010", "华商价值精选混合",             "4.93",  "89.49", "3.09", "0.1523", 10),
    @(4, "005741", "南方君信灵活配置混合A",        "3.69",  "76.52", "3.93", "0.1450", 1),
    @(5, "000029", "富国宏观策略灵活配置混合",     "5.97",  "90.27", "2.15", "0.1284", 7),
    @(6, "001364", "大成景润灵活配置混合",         "4.72",  "26.73", "2.08", "0.0982", 4),
    @(7, "002383", "大成趋势回报灵活配置混合",     "4.64",  "24.42", "1.61", "0.0747", 5),
    @(8, "003373", "大成景禄灵活配置混合A",        "3.83",  "28.91", "1.71", "0.0655", 10),
    @(9, "003374", "大成景禄灵活配置混合C",        "1.93",  "28.91", "1.71", "0.0330", 10),
    @(10, "003147", "大成动态量化配置策略混合",     "1.50",  "27.36", "2.07", "0.0310", 5),
    @(11, "630006", "华商产业升级混合",             "0.98",  "87.95", "3.09", "0.0303", 10),
    @(12, "005357", "富国国企改革灵活配置混合",     "1.13",  "87.21", "2.05", "0.0232", 9),
    @(13, "005953", "人保转型新动力灵活配置混合A",  "0.91",  "84.38", "2.15", "0.0196", 7),
    @(14, "005161", "华商上游产业股票",             "0.36",  "89.02", "3.66", "0.0132", 4),
    @(15, "006354", "国泰民裕进取灵活配置混合",     "0.52",  "80.18", "1.76", "0.0092", 9),
    @(16, "006700", "红土创新稳健混合A",            "0.74",  "27.06", "0.54", "0.0040", 8),
    @(17, "002908", "富国睿利定期开放混合",         "0.36",  "28.52", "0.85", "0.0031", 10),
    @(18, "005954", "人保转型新动力灵活配置混合C",  "0.11",  "84.38", "2.15", "0.0024", 7),
    @(19, "006701", "红土创新稳健混合C",            "0.35",  "27.06", "0.54", "0.0019", 8),
    @(20, "010150", "南方君信灵活配置混合C",        "0.02",  "76.52", "3.93", "0.0008", 1)
)

foreach ($row in $q1Rows) {
    $r = [int]$row[0] + 2
    $q1.Range("A$r").Value = [int]$row[0]
    $q1.Range("B$r").Value = "'" + $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = "'" + $row[3]
    $q1.Range("E$r").Value = "'" + $row[4]
    $q1.Range("F$r").Value = "'" + $row[5]
    $q1.Range("G$r").Value = "'" + $row[6]
    $q1.Range("H$r").Value = [int]$row[7]
}

# ---------------------------------------------------------------------------
# Step 2: refresh the (cloned) "总计" sheet with the updated totals table,
# now including the new 2022-Q1 row.
# ---------------------------------------------------------------------------

# The clone already has 6 data rows (A2:A6) styled like the original; extend
# that same index-column style one more row, down to A7.
$total.Range("A2").Copy()
$total.Range("A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 21, 2.6),
    @(1, "2021-Q4", 10, 0.38),
    @(2, "2021-Q3", 10, 0.4),
    @(3, "2021-Q2", 2, 0.1),
    @(4, "2021-Q1", 5, 0.43),
    @(5, "2020-Q4", 2, 0.31)
)

foreach ($row in $totalRows) {
    $r = [int]$row[0] + 2
    $total.Range("A$r").Value = [int]$row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = [int]$row[2]
    $total.Range("D$r").Value = [double]$row[3]
}

# Restore the originally-active sheet/tab selection (unchanged by this edit).
$wb.Worksheets.Item(1).Activate()
